$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Participant 1")

# Row 2 - first participant response
$ws.Range("B2").Value = "mission0"
$ws.Range("C2").Value = "agentPolicy0"
$ws.Range("D2").Value = "yes"

# Update header: "Time" -> "Time (minutes)"
$ws.Range("E1").Value = "Time (minutes)"

$ws.Range("E2").Value = 3.48
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = "yes"

# Row 3
$ws.Range("B3").Value = "mission13"
$ws.Range("C3").Value = "agentPolicy1"
$ws.Range("D3").Value = "no"
$ws.Range("E3").Value = 3.26
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = "no"

# Row 4
$ws.Range("B4").Value = "mission26"
$ws.Range("C4").Value = "agentPolicy2"
$ws.Range("D4").Value = "yes"
$ws.Range("E4").Value = 4.33
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = "no"

# Notes, added out of row order: row4 note first, then row2 note
$ws.Range("I4").Value = "Confident that optimality score of answer is >= 0.95"
$ws.Range("I2").Value = "Confident that optimality score of answer is >=0.95"

# Column width adjustments (bestFit-style autosize to content)
$ws.Columns.Item(5).ColumnWidth = 13.0
$ws.Columns.Item(9).ColumnWidth = 43.833333333333336

# Selection moves to I2 after data entry
$ws.Range("I2").Select()
